$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column header in H1, matching the formatting of the
# existing header row (bold font, border, centered alignment - style index 1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add data values for the new Save column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
